# Added login functionality. Began work on sign up page and functionality.
#
# The roster sheet is restructured:
#   - The old single "Name" column is split into "First Name" / "Last Name"
#     leading columns.
#   - Two new trailing columns are introduced: "Type" and "Password".
#   - Two new people are added to the roster (Nino Hana, Saki Yaki).
#   - The Email column (now column D) keeps the mailto: hyperlinks +
#     "Hyperlink" cell style that used to live on column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlinks (they lived on C2:C4) before the columns holding
# the email addresses move from C to D.
$ws.Range("C2").Hyperlinks.Delete()

# C2:C4 used to be the (styled) Email column; it becomes the plain
# Student Number column, so strip the inherited Hyperlink style.
$ws.Range("C2:C4").Style = "Normal"

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last Name"
$ws.Range("C1").Value = "Student Number"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Type"
$ws.Range("F1").Value = "Password"

# --- Data rows ---------------------------------------------------------------
# Row 2: Daniel Kolocka
$ws.Range("A2").Value = "Daniel"
$ws.Range("B2").Value = "Kolocka"
$ws.Range("C2").Value = 301333797
$ws.Range("D2").Value = "dkolocka@sfu.ca"
$ws.Range("E2").Value = "STUDENT"
$ws.Range("F2").Value = "password12"

# Row 3: Donald Trump
$ws.Range("A3").Value = "Donald"
$ws.Range("B3").Value = "Trump"
$ws.Range("C3").Value = 301333333
$ws.Range("D3").Value = "dtrump@usa.com"
$ws.Range("E3").Value = "PROFESSOR"
$ws.Range("F3").Value = "password123"

# Row 4: Mike Pence
$ws.Range("A4").Value = "Mike"
$ws.Range("B4").Value = "Pence"
$ws.Range("C4").Value = 3012222222
$ws.Range("D4").Value = "mpence@usa.com"
$ws.Range("E4").Value = "TA"
$ws.Range("F4").Value = "iamthesupreme"

# Row 5: Nino Hana (new)
$ws.Range("A5").Value = "Nino"
$ws.Range("B5").Value = "Hana"
$ws.Range("C5").Value = 1514131211
$ws.Range("D5").Value = "nhana@eorzea.ca"
$ws.Range("E5").Value = "PROFESSOR"
$ws.Range("F5").Value = "bloodforthelilies"

# Row 6: Saki Yaki (new)
$ws.Range("A6").Value = "Saki"
$ws.Range("B6").Value = "Yaki"
$ws.Range("C6").Value = 384920312
$ws.Range("D6").Value = "beefy@sfu.ca"
$ws.Range("E6").Value = "STUDENT"
$ws.Range("F6").Value = "100percentorangejuice"

# --- Email column hyperlinks ------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:dkolocka@sfu.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:dtrump@usa.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:mpence@usa.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:nhana@eorzea.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:beefy@sfu.ca") | Out-Null

# Adding a hyperlink re-writes the cell text to the target address, so put
# the display text (the email itself) back and restore the shared
# "Hyperlink" cell style (instead of the ad-hoc copy Add() applies).
$ws.Range("D2").Value = "dkolocka@sfu.ca"
$ws.Range("D3").Value = "dtrump@usa.com"
$ws.Range("D4").Value = "mpence@usa.com"
$ws.Range("D5").Value = "nhana@eorzea.ca"
$ws.Range("D6").Value = "beefy@sfu.ca"
$ws.Range("D2:D6").Style = "Hyperlink"

# --- Column widths -----------------------------------------------------------
# Target (OOXML) widths are 23.140625 / 30 / 15.85546875 / 17.5703125 /
# 14.85546875 / 29.7109375 characters. This host's ColumnWidth setter always
# re-quantizes the stored width to the nearest 1/6th of a character plus a
# fixed 5px/6 padding bump, so the literals below are the pre-images that
# land closest (nearest 1/6-character grid point) to those exact targets.
$ws.Columns.Item(1).ColumnWidth = 22.333333333333332
$ws.Columns.Item(2).ColumnWidth = 29.166666666666668
$ws.Columns.Item(3).ColumnWidth = 15.0
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.0
$ws.Columns.Item(6).ColumnWidth = 28.833333333333332

# --- Selection ----------------------------------------------------------------
$ws.Range("B13").Select()
